# edit.ps1 - reproduces the commit:
#   1. Change the table style ("Medium Style 2 - Accent 1" etc.) applied to
#      the table on slide 16 to a different built-in table style.
#   2. Re-colour the deck's master theme so the slides use the stock
#      PowerPoint "Office Theme" colour palette instead of the "Integral"
#      palette that shipped with the template (the underlying commit swaps
#      the whole ppt/theme/theme1.xml <-> ppt/theme/theme2.xml contents;
#      the closest reachable equivalent through the PowerPoint object model
#      is re-pointing every theme colour slot on the slide master to the
#      Office Theme RGB values).

$p = $ppt.ActivePresentation

# --- 1. Table style id on slide 16 ------------------------------------
$tableSlide = $p.Slides.Item(16)
$tableShape = $tableSlide.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{03C8FE8F-8E8B-47D5-B24C-B0E77EDCD199}")

# --- 2. Swap the active theme palette to the Office Theme colours -----
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

$colors.Item(1).RGB  = 0        # dk1       000000
$colors.Item(2).RGB  = 16777215 # lt1       FFFFFF
$colors.Item(3).RGB  = 6968388  # dk2       44546A
$colors.Item(4).RGB  = 15132391 # lt2       E7E6E6
$colors.Item(5).RGB  = 13998939 # accent1   5B9BD5
$colors.Item(6).RGB  = 3243501  # accent2   ED7D31
$colors.Item(7).RGB  = 10855845 # accent3   A5A5A5
$colors.Item(8).RGB  = 49407    # accent4   FFC000
$colors.Item(9).RGB  = 12874308 # accent5   4472C4
$colors.Item(10).RGB = 4697456  # accent6   70AD47
$colors.Item(11).RGB = 12673797 # hlink     0563C1
$colors.Item(12).RGB = 7491477  # folHlink  954F72
